# test_positions_onliner.xlsx — "worked a bit on the 'their_price' column in final output"
#
# - adds a new sheet "Лист4" (after "Лист2") containing a copy of the first
#   two data rows (same rows already duplicated onto "Лист1"/"Лист2"),
# - makes the new sheet the active tab,
# - narrows the selection remembered on "Лист1" to the first two rows,
# - leaves a C9 selection on the new sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet at the end of the tab strip. A freshly-added sheet is
# assigned the next sequential internal sheetId; adding (and discarding) a
# throwaway sheet first reproduces the sheetId="4" (rather than "3") that
# the target workbook ended up with.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$scratch = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)

$scratch.Delete()

# Re-fetch the new sheet by position (fresh COM reference) before activating
# it, so the workbook's remembered active-tab index is recomputed cleanly
# after the scratch sheet was removed.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Лист4"

# Populate it with the same first two rows already present on Лист1/Лист2
# (copied in two pieces so the never-used column D doesn't get touched,
# matching the source rows exactly).
$ws1.Range("A1:C2").Copy($newSheet.Range("A1:C2"))
$ws1.Range("E1:S2").Copy($newSheet.Range("E1:S2"))

# Лист1 keeps a remembered selection, just shrunk to the first two rows.
$ws1.Activate()
$ws1.Range("A1:XFD2").Select()

# The new sheet becomes the active tab, with a C9 selection.
$newSheet.Activate()
$newSheet.Range("C9").Select()
